# إضافة حدث جديد في Card19 by admin at 2025-12-08 09:22:26
#
# Card19 (row 13) previously ended with a partially-filled maintenance
# record. The edit:
#   1. Back-fills the empty columns of the existing row 13 with the
#      literal placeholder "nan" (matching how every other row in this
#      sheet represents "no value").
#   2. Appends a brand-new row 14 for the latest service event.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card19")

# --- 1. Back-fill row 13's empty cells with the "nan" placeholder ------
foreach ($col in @("B","C","D","E","F","G","H","I","J","K","M")) {
    $ws.Range($col + "13").Value = "nan"
}
# L13, N13 and O13 already hold their final values and are left untouched.

# --- 2. Append the new event as row 14 ---------------------------------
# card number - force text so it isn't stored as the number 19
$ws.Range("A14").Value = "'19"

# Date / correction / serviced-by for the new half-year service event
$ws.Range("L14").Value = "10\12\2024"
$ws.Range("N14").Value = "تم عمل صيانه نصف سنويه"
$ws.Range("O14").Value = "تيم العمل"
